$wb = $excel.ActiveWorkbook

# ---- "bb fuel upstream" sheet (new, inserted after "bb heat") ----
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws4.Name = "bb fuel upstream"
$ws4.Range("A1").Value = "meta-process"
$ws4.Range("B1").Value = "KnownQty"
$ws4.Range("C1").Value = "k_QtyFrom"
$ws4.Range("D1").Value = "UnknownQty"
$ws4.Range("E1").Value = "u_QtyTo"
$ws4.Range("F1").Value = "Calculation"
$ws4.Range("G1").Value = "Variable"
$ws4.Range("H1").Value = "2nd Known Substance"
$ws4.Range("I1").Value = "2Qty Origin"
$ws4.Range("J1").Value = "meta-notes"
$ws4.Range("B2").Value = "fuel"
$ws4.Range("C2").Value = "outflows"
$ws4.Range("D2").Value = "energy in fuel"
$ws4.Range("E2").Value = "outflows"
$ws4.Range("F2").Value = "lookup ratio-fuels"
$ws4.Range("G2").Value = "LHV"
$ws4.Range("B3").Value = "fuel"
$ws4.Range("C3").Value = "outflows"
$ws4.Range("D3").Value = "fuel"
$ws4.Range("E3").Value = "inflow"
$ws4.Range("F3").Value = "returnvalue"
$ws4.Range("B4").Value = "fuel"
$ws4.Range("C4").Value = "outflows"
$ws4.Range("D4").Value = "CO2__emitted"
$ws4.Range("E4").Value = "outflows"
$ws4.Range("F4").Value = "lookup ratio-fuels"
$ws4.Range("G4").Value = "upstream CO2"
$ws4.Range("B5").Value = "CO2__emitted"
$ws4.Range("C5").Value = "outflows"
$ws4.Range("D5").Value = "carbon and oxygen"
$ws4.Range("E5").Value = "inflows"
$ws4.Range("F5").Value = "returnvalue"
$ws4.Rows.Item(1).Font.Size = 11
$ws4.Rows.Item(1).Font.Bold = $true
$ws4.Columns.Item(2).AutoFit() | Out-Null
$ws4.Columns.Item(3).AutoFit() | Out-Null
$ws4.Columns.Item(4).AutoFit() | Out-Null
$ws4.Columns.Item(5).AutoFit() | Out-Null
$ws4.Columns.Item(6).AutoFit() | Out-Null
$ws4.Range("B5").Select() | Out-Null

# ---- "bb biofuel upstream" sheet (new, inserted after the previous one, ends up active) ----
$ws5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws5.Name = "bb biofuel upstream"
$ws5.Range("A1").Value = "meta-process"
$ws5.Range("B1").Value = "KnownQty"
$ws5.Range("C1").Value = "k_QtyFrom"
$ws5.Range("D1").Value = "UnknownQty"
$ws5.Range("E1").Value = "u_QtyTo"
$ws5.Range("F1").Value = "Calculation"
$ws5.Range("G1").Value = "Variable"
$ws5.Range("H1").Value = "2nd Known Substance"
$ws5.Range("I1").Value = "2Qty Origin"
$ws5.Range("J1").Value = "meta-notes"
$ws5.Range("B2").Value = "biofuel"
$ws5.Range("C2").Value = "outflows"
$ws5.Range("D2").Value = "energy in biofuel"
$ws5.Range("E2").Value = "outflows"
$ws5.Range("F2").Value = "lookup ratio"
$ws5.Range("G2").Value = "LHV"
$ws5.Range("B3").Value = "biofuel"
$ws5.Range("C3").Value = "outflows"
$ws5.Range("D3").Value = "fresh biomass"
$ws5.Range("E3").Value = "inflow"
$ws5.Range("F3").Value = "lookup ratio-fuels"
$ws5.Range("G3").Value = "fresh biomass ratio"
$ws5.Range("B4").Value = "fresh biomass"
$ws5.Range("C4").Value = "inflow"
$ws5.Range("D4").Value = "biomass losses"
$ws5.Range("E4").Value = "outflows"
$ws5.Range("F4").Value = "subtraction"
$ws5.Range("H4").Value = "biofuel"
$ws5.Range("I4").Value = "outflows"
$ws5.Range("B5").Value = "fresh biomass"
$ws5.Range("C5").Value = "inflow"
$ws5.Range("D5").Value = "biofuel__biomass"
$ws5.Range("E5").Value = "temp"
$ws5.Range("F5").Value = "returnvalue"
$ws5.Range("B6").Value = "biofuel__biomass"
$ws5.Range("C6").Value = "temp"
$ws5.Range("D6").Value = "CO2__removed from atmosphere"
$ws5.Range("E6").Value = "inflows"
$ws5.Range("F6").Value = "lookup ratio-fuels"
$ws5.Range("G6").Value = "biomass CO2 absorption"
$ws5.Range("B7").Value = "biofuel"
$ws5.Range("C7").Value = "outflows"
$ws5.Range("D7").Value = "CO2__emitted"
$ws5.Range("E7").Value = "outflows"
$ws5.Range("F7").Value = "lookup ratio-fuels"
$ws5.Range("G7").Value = "upstream CO2"
$ws5.Range("B8").Value = "CO2__removed from atmosphere"
$ws5.Range("C8").Value = "inflows"
$ws5.Range("D8").Value = "CONSUMED CO2 removals net emissions"
$ws5.Range("E8").Value = "outflows"
$ws5.Range("F8").Value = "subtraction"
$ws5.Range("H8").Value = "CO2__emitted"
$ws5.Range("I8").Value = "outflows"
$ws5.Rows.Item(1).Font.Size = 11
$ws5.Rows.Item(1).Font.Bold = $true
$ws5.Range("G3").WrapText = $true
$ws5.Columns.Item(1).AutoFit() | Out-Null
$ws5.Columns.Item(2).AutoFit() | Out-Null
$ws5.Columns.Item(3).AutoFit() | Out-Null
$ws5.Columns.Item(4).AutoFit() | Out-Null
$ws5.Columns.Item(5).AutoFit() | Out-Null
$ws5.Columns.Item(6).AutoFit() | Out-Null
$ws5.Columns.Item(7).AutoFit() | Out-Null

# ---- restore the "bb heat" selection to a plain single cell & move the active tab ----
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("H29").Select() | Out-Null

$ws5.Activate() | Out-Null
$ws5.Range("E19").Select() | Out-Null

Write-Output "done"
